$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D to Text so numeric-looking price strings
# (e.g. "1.001", "306.24") are stored as text, matching the source data,
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.113.08'
$ws.Range("E2").Value = '  -2.53%  '
$ws.Range("D3").Value = '1.865.13'
$ws.Range("E3").Value = '  -2.37%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").Value = '306.24'
$ws.Range("E5").Value = '  -2.09%  '
$ws.Range("E6").Value = '  +0.14%  '
$ws.Range("D7").Value = '0.5153'
$ws.Range("E7").Value = '  -2.11%  '
$ws.Range("E8").Value = '  -0.62%  '
$ws.Range("D9").Value = '0.07155'
$ws.Range("E9").Value = '  -1.30%  '
$ws.Range("D10").Value = '0.8890'
$ws.Range("E10").Value = '  -1.97%  '
$ws.Range("D11").Value = '20.70'
$ws.Range("E11").Value = '  -3.08%  '
$ws.Range("D12").Value = '0.07554'
$ws.Range("E12").Value = '  -1.33%  '
$ws.Range("D13").Value = '1.835.23'
$ws.Range("E13").Value = '  -3.84%  '
$ws.Range("D14").Value = '5.309'
$ws.Range("E14").Value = '  -2.72%  '
$ws.Range("D15").Value = '89.73'
$ws.Range("E15").Value = '  -2.63%  '
$ws.Range("E16").Value = '  +0.28%  '
$ws.Range("D17").Value = '0.000008466'
$ws.Range("E17").Value = '  -3.02%  '
$ws.Range("D18").Value = '14.04'
$ws.Range("E18").Value = '  -3.47%  '
$ws.Range("D19").Value = '1.000'
$ws.Range("E19").Value = '  +0.02%  '
$ws.Range("D20").Value = '27.143.22'
$ws.Range("E20").Value = '  -2.56%  '
$ws.Range("D21").Value = '5.014'
$ws.Range("E21").Value = '  -2.79%  '
$ws.Range("D22").Value = '2.088.60'
$ws.Range("E22").Value = '  -5.37%  '
$ws.Range("E23").Value = '  -3.45%  '
$ws.Range("E24").Value = '  -3.13%  '
$ws.Range("D25").Value = '1.839'
$ws.Range("E25").Value = '  -1.67%  '
$ws.Range("D26").Value = '145.77'
$ws.Range("E26").Value = '  -5.26%  '
$ws.Range("D27").Value = '17.94'
$ws.Range("E27").Value = '  -2.32%  '
$ws.Range("D28").Value = '2.089'
$ws.Range("E28").Value = '  -3.85%  '
$ws.Range("D29").Value = '112.73'
$ws.Range("E29").Value = '  -1.93%  '
$ws.Range("D30").Value = '4.665'
$ws.Range("E30").Value = '  -4.06%  '
$ws.Range("D31").Value = '4.655'
$ws.Range("E31").Value = '  -4.33%  '
$ws.Range("D32").Value = '0.09153'
$ws.Range("E32").Value = '  +0.57%  '
$ws.Range("D33").Value = '0.05098'
$ws.Range("E33").Value = '  -3.54%  '
$ws.Range("D34").Value = '3.071'
$ws.Range("E34").Value = '  -3.54%  '
$ws.Range("D35").Value = '1.156'
$ws.Range("E35").Value = '  -6.56%  '
$ws.Range("D36").Value = '0.7241'
$ws.Range("E36").Value = '  -7.21%  '
$ws.Range("E37").Value = '  -2.89%  '
$ws.Range("D38").Value = '3.085'
$ws.Range("E38").Value = '  +0.30%  '
$ws.Range("D39").Value = '2.495'
$ws.Range("E39").Value = '  -4.64%  '
$ws.Range("D40").Value = '1.075'
$ws.Range("E40").Value = '  -1.72%  '
$ws.Range("D41").Value = '0.5276'
$ws.Range("E41").Value = '  -5.93%  '
$ws.Range("D42").Value = '6.459'
$ws.Range("E42").Value = '  -3.95%  '
$ws.Range("D43").Value = '116.02'
$ws.Range("E43").Value = '  +0.06%  '
$ws.Range("D44").Value = '8.277'
$ws.Range("E44").Value = '  -3.56%  '
$ws.Range("E45").Value = '  -3.82%  '
$ws.Range("E46").Value = '  +0.13%  '
$ws.Range("D47").Value = '0.4616'
$ws.Range("E47").Value = '  -4.34%  '
$ws.Range("D48").Value = '9.946'
$ws.Range("E48").Value = '  -5.45%  '
$ws.Range("D49").Value = '1.563'
$ws.Range("E49").Value = '  -3.64%  '
$ws.Range("D50").Value = '36.52'
$ws.Range("E50").Value = '  -1.46%  '
$ws.Range("D51").Value = '63.35'
$ws.Range("E51").Value = '  -5.53%  '

# Restore the original (default/Normal) cell style now that the text
# values are safely stored, so no stray number-format style lingers.
$ws.Range("D2:D51").Style = "Normal"
